{"js": "// The document contains a series of paragraphs shaped like\n// \"Label: value\" (e.g. \"Name: az18042003an\", \"Phone: 0332360580\", ...).\n// The edit clears out each value, leaving just \"Label: \" (colon + single\n// trailing space) behind.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = paragraph.text;\n  const sepIndex = text.indexOf(\": \");\n  if (sepIndex === -1) {\n    continue; // not a \"Label: value\" paragraph, leave untouched\n  }\n  const label = text.substring(0, sepIndex + 2); // keep \"Label: \"\n  if (label === text) {\n    continue; // already has no value, nothing to do\n  }\n  paragraph.insertText(label, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains a series of paragraphs shaped like\n# \"Label: value\" (e.g. \"Name: az18042003an\", \"Phone: 0332360580\", ...).\n# The edit clears out each value, leaving just \"Label: \" (colon + single\n# trailing space) behind.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $t = $r.Text\n    $idx = $t.IndexOf(\": \")\n    if ($idx -ge 0) {\n        $label = $t.Substring(0, $idx + 2)\n        if ($label -ne $t) {\n            $r.Text = $label\n        }\n    }\n}\n"}
